$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1/IF (same style as existing header H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# New data cells for rows 2 and 3
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 7
$ws.Range("J3").Value = 7
